# SC Charge List excel state removed, name changes
#
# The "State" column (column B) is removed from the SC Charges List
# template. All columns to its right (Product, Category, Capacity,
# Service Category, Vendor Basic Charge, Vendor Tax, Vendor Total,
# Customer Total Rs., Serial Number Mandatory) shift one column to the
# left, taking their formatting/styles with them. Deleting the column
# also drops the now-unused "State" / "{sc:state}" shared-string
# entries automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "State" column (B) entirely - content, styles and column
# widths for C:K shift left into B:J.
$ws.Columns("B").Delete()

# View-state tweaks that came along with this edit: zoom to 70% and
# move the selection to C10 (previously the view was scrolled right
# with topLeftCell F1 and the selection sat at I8).
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("C10").Select()
